$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-16 from 45243 (2023-11-13)
# to 45244 (2023-11-14), matching the daily automatic refresh.
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
